$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7:C7").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$ws.Range("B8").Value = 44232
$ws.Range("C8").Value = "11h40"
$ws.Range("D8").Value = "12h20"
$ws.Range("E8").Value = "40min"
$ws.Range("F8").Value = "Projet WEB annonces"
$ws.Range("G8").Value = "Consseption des USE CASE + Scénarios"

$ws.Range("F11").Select()
